$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01142538604652072
$ws.Range("C2").Value = 0.2150658939635976

$ws.Range("B3").Value = 0.07749699958783268
$ws.Range("C3").Value = 0.139516648354486

$ws.Range("B4").Value = 0.7293133251922297
$ws.Range("C4").Value = 0.1972748460205353

$ws.Range("B5").Value = 0.9778767664683011
$ws.Range("C5").Value = 0.3784942077030671

$ws.Range("B6").Value = 0.8939181292189049
$ws.Range("C6").Value = 0.6169331304427405

$ws.Range("B7").Value = 0.7911348074035913
$ws.Range("C7").Value = 0.1119329866628683

$ws.Range("B8").Value = 0.005086978077888489
$ws.Range("C8").Value = 0.2117326927185059
